$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting
# (values such as "280.10" or "1.00" must not be coerced to numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.140.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.377.25"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +6.79%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.43"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +9.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.84"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -6.67%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.651"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +8.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.33"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.98%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.56"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.97"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +13.63%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.737.64"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +6.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.392.13"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +7.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.123.64"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.84"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +8.82%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.67"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "280.10"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +17.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.38"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.58"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +8.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.61"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.13"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +7.55%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "176.57"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.82"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.60%  "
$ws.Range("B31").Value = "WEMIXToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.19"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.88%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.15"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0921"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.85"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.133"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.85"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.16"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0363"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.54%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.82"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +17.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.57"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +20.49%  "
$ws.Range("B42").Value = "BitcoinSV"
$ws.Range("C42").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.84"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +69.11%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "123.33"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +21.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "68.95"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.54%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.42"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.44"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +11.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.58"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.37%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.23%  "
